$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F5").Value = -8
$ws.Range("F11").Value = -6
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = -5
$ws.Range("F17").Value = -4
$ws.Range("F20").Value = -3
$ws.Range("F32").Value = 7
$ws.Range("F39").Value = -1
$ws.Range("F40").Value = 2
